$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "*maa://24633 (56.17), *maa://30515 (70.48), maa://39402 (93.15), *maa://34787 (73.42), ***maa://20792 (11.93), ***maa://29083 (27.78)"
$ws.Range("AB2").Value = "maa://21246 (91.47), maa://36684 (95.24), ***maa://22731 (6.25)"
$ws.Range("X3").Value = "maa://27396 (83.94), maa://27484 (96.75), maa://27480 (83.33)"
$ws.Range("AB3").Value = "maa://24390 (95.0)"
$ws.Range("T4").Value = "maa://32509 (95.87), maa://27295 (87.18), maa://22754 (90.41), *maa://31008 (78.57), *maa://21746 (55.81)"
$ws.Range("X4").Value = "**maa://32495 (48.54), ***maa://31785 (22.22), maa://43217 (90.59), ***maa://36683 (29.79)"
$ws.Range("D6").Value = "maa://42407 (95.38)"
$ws.Range("A8").Value = "更新日期：2025.04.13 13:23:50"
$ws.Range("AB9").Value = "maa://28711 (87.3), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), *maa://45044 (66.67), maa://40166 (96.67)"
$ws.Range("D10").Value = "***maa://25695 (18.65), ***maa://39951 (12.7), ***maa://34206 (22.22), ***maa://39243 (25.0), *maa://45271 (54.17)"
$ws.Range("T11").Value = "maa://22747 (92.73), maa://22501 (97.92), maa://45521 (88.0)"
$ws.Range("AB11").Value = "maa://29912 (97.4), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("D12").Value = "maa://30766 (89.29), *maa://36678 (77.78)"
$ws.Range("D14").Value = "maa://30764 (88.33)"
$ws.Range("T15").Value = "maa://23892 (96.39)"
$ws.Range("P16").Value = "maa://28504 (91.94)"
$ws.Range("X16").Value = "maa://28501 (98.17), maa://28051 (96.0)"
$ws.Range("D17").Value = "maa://21624 (85.37)"
$ws.Range("T17").Value = "*maa://42324 (51.22)"
$ws.Range("L18").Value = "maa://22466 (91.3), *maa://22732 (52.04)"
$ws.Range("AB18").Value = "maa://24393 (98.0)"
$ws.Range("D20").Value = "maa://21432 (90.81), maa://25198 (93.81), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("L20").Value = "maa://41331 (86.11)"
$ws.Range("P20").Value = "maa://37442 (95.74)"
$ws.Range("X20").Value = "maa://49976 (91.3), maa://50085 (100.0)"
$ws.Range("AF22").Value = "maa://29658 (94.0)"
$ws.Range("L23").Value = "maa://39756 (95.87), maa://39875 (94.59)"
$ws.Range("X23").Value = "*maa://28503 (69.77)"
$ws.Range("AB23").Value = "maa://29652 (97.78)"
$ws.Range("X24").Value = "maa://29988 (83.4), maa://23504 (93.48), **maa://22892 (40.27), *maa://25141 (77.27), *maa://36663 (78.31), ***maa://22815 (23.08)"
$ws.Range("D25").Value = "maa://29753 (95.32)"
$ws.Range("P25").Value = "maa://24382 (93.75)"
$ws.Range("AB25").Value = "maa://31215 (88.62), maa://24516 (80.22), maa://26001 (87.5)"
$ws.Range("AF25").Value = "maa://20108 (96.4), maa://24621 (96.97), maa://36676 (97.06), maa://22771 (85.71), *maa://37772 (66.67)"
$ws.Range("D26").Value = "maa://41802 (92.0)"
$ws.Range("X26").Value = "maa://24389 (96.77)"
$ws.Range("AB26").Value = "maa://42235 (95.04)"
$ws.Range("AF26").Value = "maa://30511 (82.22), *maa://29760 (56.25)"
$ws.Range("L27").Value = "maa://28071 (90.91)"
$ws.Range("AF28").Value = "maa://36660 (92.49), *maa://36701 (66.67)"
$ws.Range("D29").Value = "maa://31694 (98.18)"
$ws.Range("D30").Value = "maa://45792 (90.91)"
$ws.Range("L31").Value = "maa://35926 (93.49), maa://36258 (85.16), *maa://43904 (75.0)"
$ws.Range("T31").Value = "maa://30711 (96.67), maa://30768 (100.0)"
$ws.Range("L32").Value = "maa://28065 (95.83)"
$ws.Range("AF32").Value = "maa://42408 (85.71)"
$ws.Range("P34").Value = "maa://48817 (92.59)"
$ws.Range("AF34").Value = "*maa://32650 (77.27)"
$ws.Range("L35").Value = "maa://41296 (96.74)"
$ws.Range("AF35").Value = "maa://39479 (90.48)"
$ws.Range("P37").Value = "maa://21280 (89.69), *maa://21239 (69.23)"
$ws.Range("T38").Value = "maa://30713 (97.06)"
$ws.Range("P41").Value = "**maa://35616 (40.0), maa://43177 (92.0)"
$ws.Range("T44").Value = "maa://39366 (89.47)"
$ws.Range("P45").Value = "*maa://36237 (72.22)"
$ws.Range("H46").Value = "maa://35931 (92.12), maa://43901 (94.29)"
$ws.Range("P49").Value = "*maa://39643 (65.71)"
$ws.Range("H58").Value = "*maa://37964 (57.78)"
$ws.Range("H60").Value = "*maa://40438 (72.06)"
